$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rows 3, 4, 5 got re-sorted (the match data for columns F..V rotates
#    one row down, wrapping from row 5 back to row 3). Columns A..E
#    (Indice / pais / torneio / temporada / data_partida) stay put.
#    Capture the "before" F..V values for each of the three rows first,
#    then write them back in rotated order so earlier writes don't
#    clobber values we still need to read.
# ---------------------------------------------------------------------
$cols = 6..22   # F..V

$row3 = @{}
$row4 = @{}
$row5 = @{}
foreach ($c in $cols) {
    $row3[$c] = $ws.Cells.Item(3, $c).Value()
    $row4[$c] = $ws.Cells.Item(4, $c).Value()
    $row5[$c] = $ws.Cells.Item(5, $c).Value()
}

foreach ($c in $cols) {
    $ws.Cells.Item(3, $c).Value = $row5[$c]
    $ws.Cells.Item(4, $c).Value = $row3[$c]
    $ws.Cells.Item(5, $c).Value = $row4[$c]
}

# ---------------------------------------------------------------------
# 2) Two new match rows were appended at the bottom (62 and 63).
#    Insert them by duplicating the last existing row (61) so the new
#    rows inherit the same look (border/alignment/bold on column A,
#    datetime format on column E), then overwrite with the real data.
# ---------------------------------------------------------------------
function New-MatchRow($rowNum, $data) {
    $ws.Rows(61).Copy()
    $ws.Rows($rowNum).Insert()
    # The engine's row-insert drops the border on the copied style for
    # column A; put it back so the style matches the original (bold,
    # centered, thin box border) instead of minting a near-duplicate one.
    $ws.Cells.Item($rowNum, 1).Borders.LineStyle = 1
    $ws.Cells.Item($rowNum, 1).Borders.Weight = 2

    $ws.Cells.Item($rowNum, 1).Value = $data.A
    $ws.Cells.Item($rowNum, 2).Value = $data.B
    $ws.Cells.Item($rowNum, 3).Value = $data.C
    $ws.Cells.Item($rowNum, 4).Value = $data.D
    $ws.Cells.Item($rowNum, 5).Value = $data.E
    $ws.Cells.Item($rowNum, 6).Value = $data.F
    $ws.Cells.Item($rowNum, 7).Value = $data.G
    $ws.Cells.Item($rowNum, 8).Value = $data.H
    $ws.Cells.Item($rowNum, 9).Value = $data.I
    $ws.Cells.Item($rowNum, 10).Value = $data.J
    $ws.Cells.Item($rowNum, 11).Value = $data.K
    $ws.Cells.Item($rowNum, 12).Value = $data.L
    $ws.Cells.Item($rowNum, 13).Value = $data.M
    $ws.Cells.Item($rowNum, 14).Value = $data.N
    $ws.Cells.Item($rowNum, 15).Value = $data.O
    $ws.Cells.Item($rowNum, 16).Value = $data.P
    $ws.Cells.Item($rowNum, 17).Value = $data.Q
    $ws.Cells.Item($rowNum, 18).Value = $data.R
    $ws.Cells.Item($rowNum, 19).Value = $data.S
    $ws.Cells.Item($rowNum, 20).Value = $data.T
    $ws.Cells.Item($rowNum, 21).Value = $data.U
    $ws.Cells.Item($rowNum, 22).Value = $data.V
}

$data62 = @{
    A = 61; B = "iran"; C = "persian-gulf-pro-league"; D = "2023-2024"
    E = 45233.52083333334
    F = "Havadar SC"; G = 0; H = "Paykan"; I = 0
    J = 2.17; K = "02/11/2023 00:42"
    L = 2.15; M = "03/11/2023 12:28"
    N = 2.65; O = "02/11/2023 00:42"
    P = 2.41; Q = "03/11/2023 12:28"
    R = 3.62; S = "02/11/2023 00:42"
    T = 4.31; U = "03/11/2023 12:28"
    V = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/havadar-sc-paykan/KvfBkQWq/"
}

$data63 = @{
    A = 62; B = "iran"; C = "persian-gulf-pro-league"; D = "2023-2024"
    E = 45233.52083333334
    F = "Mes Rafsanjan"; G = 3; H = "Foolad"; I = 0
    J = 2.2; K = "02/11/2023 00:42"
    L = 2.3; M = "03/11/2023 12:29"
    N = 2.65; O = "02/11/2023 00:42"
    P = 2.48; Q = "03/11/2023 12:29"
    R = 3.55; S = "02/11/2023 00:42"
    T = 4.29; U = "03/11/2023 12:29"
    V = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/mes-rafsanjan-foolad/z7Iven9M/"
}

New-MatchRow 62 $data62
New-MatchRow 63 $data63
